$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (SCA_N)
$ws.Range("B2").Value = 0.7446105873521773
$ws.Range("C2").Value = 0.6085072646245111
$ws.Range("D2").Value = -0.57049515031306

# Row 3 (EA_N)
$ws.Range("B3").Value = 0.7122635255421865
$ws.Range("C3").Value = -0.7132032837917984
$ws.Range("D3").Value = 0.629619077900694

# Row 4 (ENSO-mei_N)
$ws.Range("B4").Value = 0.7392183375353897
$ws.Range("C4").Value = 0.714741607648563
$ws.Range("D4").Value = -0.8246507296936757

# Row 5 (NAO_N)
$ws.Range("B5").Value = 0.7469248575295034
$ws.Range("C5").Value = -0.6347582799799147
$ws.Range("D5").Value = -0.8196024610864064

# Row 6 (SCA_P)
$ws.Range("B6").Value = 0.6857279732129107
$ws.Range("C6").Value = -0.6055451143102873
$ws.Range("D6").Value = -0.7428956286176914

# Row 7 (EA_P) - B7 stays empty
$ws.Range("C7").Value = 0.5633557605361711
$ws.Range("D7").Value = 0.5146089570365168

# Row 8 (ENSO-mei_P) - D8 stays empty
$ws.Range("B8").Value = 0.5981877440692995
$ws.Range("C8").Value = 0.6691189131641496

# Row 9 (NAO_P) - B9 stays empty
$ws.Range("C9").Value = 0.6650386990886288
$ws.Range("D9").Value = 0.6590298530350931
